$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (columns C, F, G) starting at row 79
$data = @(
    @("Plain",  28.3733,             6),
    @(60079,    6.3574999999999999, 18),
    @(45096,    27.905799999999999,  3),
    @(188091,   10.9636,             7),
    @("test",   20.583500000000001,  9),
    @(253036,   13.8782,             4),
    @(42049,    13.3101,             6),
    @(35070,    9.7677999999999994, 14),
    @(163014,   10.9528,            14),
    @(124084,   26.652000000000001,  4)
)

$startRow = 79
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 6).Value = $data[$i][1]
    $ws.Cells.Item($row, 7).Value = $data[$i][2]
}

$ws.Range("E81").Select()
